$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.068.49"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.46%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.313.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.04%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "302.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.43%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.503"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.88%  "

$ws.Range("E8").Value = "  -0.04%  "

$ws.Range("E9").Value = "  +5.39%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.15"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +9.81%  "

$ws.Range("E11").Value = "  +0.99%  "

$ws.Range("E12").Value = "  +2.91%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "17.85"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +13.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.43%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.673.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.309.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.62%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.807"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.965.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.66%  "

$ws.Range("E20").Value = "  +4.30%  "

$ws.Range("E21").Value = "  +1.64%  "

$ws.Range("E22").Value = "  +2.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "236.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.49%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.21"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +13.67%  "

$ws.Range("B25").Value = "Dai"
$ws.Range("C25").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.00"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.06%  "

$ws.Range("B26").Value = "PancakeSwap"
$ws.Range("C26").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.46%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.30%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +8.61%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "34.80"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.72%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "168.61"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.55%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "9.19"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.87%  "

$ws.Range("E32").Value = "  +0.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.01"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.82%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.67"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.23%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.45"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.51%  "

$ws.Range("E36").Value = "  +2.28%  "

$ws.Range("E37").Value = "  +1.26%  "

$ws.Range("E38").Value = "  +4.42%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.85"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.17%  "

$ws.Range("E40").Value = "  +4.43%  "

$ws.Range("E41").Value = "  +1.11%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.985.29"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.52%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0289"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.80%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.30%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.25"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.52%  "

$ws.Range("E46").Value = "  +5.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.14%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "56.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +8.29%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.539.89"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.89%  "

$ws.Range("E50").Value = "  +4.17%  "

$ws.Range("B51").Value = "THORChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.56"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.37%  "
